$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a non-UK delivery scenario for the existing "sony" order ---
# Insert a new row above the current "sony" row (row 2): same product/qty,
# but billed in the UK and delivered to Japan (i.e. delivery != billing,
# and delivery is not UK).
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "sony"
$ws.Range("B2").Value = "Sony VAIO"
$ws.Range("C2").Value = "'2"
$ws.Range("D2").Value = "United Kingdom"
$ws.Range("E2").Value = "Japan"

# --- Change the "hp" order (now row 5) to a non-UK billing/delivery case ---
$ws.Range("C5").Value = "'1"
$ws.Range("D5").Value = "Japan"
$ws.Range("E5").Value = "Japan"

# --- Replace the old "mac" order (now row 6) with another non-UK delivery
#     scenario for the "ipod" product ---
$ws.Rows.Item(6).Delete()
$ws.Range("A6").Value = "ipod"
$ws.Range("B6").Value = "iPod Classic"
$ws.Range("C6").Value = "'5"
$ws.Range("D6").Value = "United Kingdom"
$ws.Range("E6").Value = "Japan"

# Leave the cursor on the next (empty) row, as if the editor had just
# finished typing in this data set.
$ws.Range("C7").Select() | Out-Null
